# Auto-generated edit script: updates calculated profit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1087.8572
$ws.Range("J12").Value = 206.25
$ws.Range("L12").Value = 206.25
$ws.Range("N12").Value = -546.25
$ws.Range("H19").Value = 291.4
$ws.Range("I19").Value = 267.14285
$ws.Range("J19").Value = 312.625
$ws.Range("K19").Value = 267.14285
$ws.Range("L19").Value = 312.625
$ws.Range("M19").Value = -92.14285000000001
$ws.Range("N19").Value = -662.625
$ws.Range("H40").Value = 1045.25
$ws.Range("I40").Value = 792.8461
$ws.Range("J40").Value = 1343.5454
$ws.Range("K40").Value = 792.8461
$ws.Range("L40").Value = 1343.5454
$ws.Range("M40").Value = -617.8461
$ws.Range("N40").Value = -1693.5454
$ws.Range("H70").Value = 1032.375
$ws.Range("I70").Value = 1286.6666
$ws.Range("J70").Value = 879.8
$ws.Range("K70").Value = 3859.9998
$ws.Range("L70").Value = 2639.4
$ws.Range("M70").Value = -3589.9998
$ws.Range("N70").Value = -3179.4
$ws.Range("H73").Value = 1032.375
$ws.Range("I73").Value = 1286.6666
$ws.Range("J73").Value = 879.8
$ws.Range("K73").Value = 3859.9998
$ws.Range("L73").Value = 2639.4
$ws.Range("M73").Value = -2923.9998
$ws.Range("N73").Value = -4511.4
$ws.Range("H101").Value = 240
$ws.Range("J101").Value = 300
$ws.Range("L101").Value = 900
$ws.Range("N101").Value = -4144
$ws.Range("H116").Value = 22732772
$ws.Range("I116").Value = 83335496
$ws.Range("J116").Value = 6749.5
$ws.Range("K116").Value = 83335496
$ws.Range("L116").Value = 6749.5
$ws.Range("M116").Value = -83332054
$ws.Range("N116").Value = -13633.5
$ws.Range("H132").Value = 4754.864
$ws.Range("I132").Value = 5217.706
$ws.Range("J132").Value = 3181.2
$ws.Range("K132").Value = 15653.118
$ws.Range("L132").Value = 9543.599999999999
$ws.Range("M132").Value = -13123.118
$ws.Range("N132").Value = -14603.6
$ws.Range("H138").Value = 2678.5403
$ws.Range("I138").Value = 1821.875
$ws.Range("J138").Value = 2871.5916
$ws.Range("K138").Value = 5465.625
$ws.Range("L138").Value = 8614.774800000001
$ws.Range("M138").Value = -325.625
$ws.Range("N138").Value = -18894.7748

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1390.3
$ws.Range("I2").Value = 1398.5714
$ws.Range("K2").Value = 1398.5714
$ws.Range("M2").Value = -1285.5714
$ws.Range("H45").Value = 2483.2258
$ws.Range("I45").Value = 2388.5
$ws.Range("K45").Value = 2388.5
$ws.Range("M45").Value = -2011.5
$ws.Range("H74").Value = 47621136
$ws.Range("I74").Value = 76923820
$ws.Range("J74").Value = 4274.875
$ws.Range("K74").Value = 76923820
$ws.Range("L74").Value = 4274.875
$ws.Range("M74").Value = -76922946
$ws.Range("N74").Value = -6022.875
$ws.Range("H77").Value = 47621136
$ws.Range("I77").Value = 76923820
$ws.Range("J77").Value = 4274.875
$ws.Range("K77").Value = 384619100
$ws.Range("L77").Value = 21374.375
$ws.Range("M77").Value = -384614732
$ws.Range("N77").Value = -30110.375
$ws.Range("H116").Value = 1390.3
$ws.Range("I116").Value = 1398.5714
$ws.Range("K116").Value = 1398.5714
$ws.Range("M116").Value = 895.4286
$ws.Range("H122").Value = 1904.8485
$ws.Range("I122").Value = 1763.2142
$ws.Range("K122").Value = 5289.642599999999
$ws.Range("M122").Value = -2839.642599999999
$ws.Range("H132").Value = 13065.889
$ws.Range("I132").Value = 1769.6285
$ws.Range("K132").Value = 5308.8855
$ws.Range("M132").Value = -2778.8855

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1390.3
$ws.Range("I3").Value = 1398.5714
$ws.Range("K3").Value = 1398.5714
$ws.Range("M3").Value = -1284.5714
$ws.Range("H20").Value = 2515.6316
$ws.Range("I20").Value = 2360.8
$ws.Range("K20").Value = 2360.8
$ws.Range("M20").Value = -2113.8
$ws.Range("H134").Value = 3909.6365
$ws.Range("I134").Value = 4123.467
$ws.Range("K134").Value = 12370.401
$ws.Range("M134").Value = -9835.400999999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4395.8604
$ws.Range("I31").Value = 2032.4286
$ws.Range("J31").Value = 6651.864
$ws.Range("K31").Value = 2032.4286
$ws.Range("L31").Value = 6651.864
$ws.Range("M31").Value = -1737.4286
$ws.Range("N31").Value = -7241.864
$ws.Range("H34").Value = 4395.8604
$ws.Range("I34").Value = 2032.4286
$ws.Range("J34").Value = 6651.864
$ws.Range("K34").Value = 2032.4286
$ws.Range("L34").Value = 6651.864
$ws.Range("M34").Value = -1830.4286
$ws.Range("N34").Value = -7055.864
$ws.Range("H94").Value = 3616.3125
$ws.Range("I94").Value = 2429.2222
$ws.Range("J94").Value = 5142.5713
$ws.Range("K94").Value = 2429.2222
$ws.Range("L94").Value = 5142.5713
$ws.Range("M94").Value = -1978.2222
$ws.Range("N94").Value = -6044.5713
$ws.Range("H105").Value = 7813268
$ws.Range("I105").Value = 12500658
$ws.Range("J105").Value = 951.6667
$ws.Range("K105").Value = 12500658
$ws.Range("L105").Value = 951.6667
$ws.Range("M105").Value = -12498911
$ws.Range("N105").Value = -4445.6667
$ws.Range("H109").Value = 76176450
$ws.Range("J109").Value = 76176450
$ws.Range("L109").Value = 76176450
$ws.Range("N109").Value = -76178530
$ws.Range("H122").Value = 1260.1904
$ws.Range("I122").Value = 946
$ws.Range("J122").Value = 2045.6666
$ws.Range("K122").Value = 2838
$ws.Range("L122").Value = 6136.9998
$ws.Range("M122").Value = -388
$ws.Range("N122").Value = -11036.9998
$ws.Range("H132").Value = 2627.2903
$ws.Range("I132").Value = 1736.6957
$ws.Range("K132").Value = 5210.0871
$ws.Range("M132").Value = -2680.0871
$ws.Range("H134").Value = 790.9091
$ws.Range("I134").Value = 587.5
$ws.Range("K134").Value = 1762.5
$ws.Range("M134").Value = 772.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 3050
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3050
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 9150
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = -11522
$ws.Range("H89").Value = 3050
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3050
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 27450
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = -39306

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 995.39026
$ws.Range("I102").Value = 769.925
$ws.Range("J102").Value = 10014
$ws.Range("K102").Value = 769.925
$ws.Range("L102").Value = 10014
$ws.Range("M102").Value = 852.075
$ws.Range("N102").Value = -13258
$ws.Range("H139").Value = 39006
$ws.Range("J139").Value = 39006
$ws.Range("L139").Value = 39006
$ws.Range("N139").Value = -49286

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4120
$ws.Range("I7").Value = 3722.2222
$ws.Range("K7").Value = 3722.2222
$ws.Range("M7").Value = -3610.2222
$ws.Range("H14").Value = 1633.3334
$ws.Range("J14").Value = 1633.3334
$ws.Range("L14").Value = 1633.3334
$ws.Range("N14").Value = -1977.3334
$ws.Range("H68").Value = 2840.8572
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2840.8572
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2840.8572
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = -4338.8572
$ws.Range("H71").Value = 2840.8572
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2840.8572
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 14204.286
$ws.Range("M71").Value = $null
$ws.Range("N71").Value = -21692.286
$ws.Range("H82").Value = 1911.1111
$ws.Range("I82").Value = 2045.7142
$ws.Range("J82").Value = 1440
$ws.Range("K82").Value = 2045.7142
$ws.Range("L82").Value = 1440
$ws.Range("M82").Value = -1684.7142
$ws.Range("N82").Value = -2162
$ws.Range("H85").Value = 1911.1111
$ws.Range("I85").Value = 2045.7142
$ws.Range("J85").Value = 1440
$ws.Range("K85").Value = 2045.7142
$ws.Range("L85").Value = 1440
$ws.Range("M85").Value = -797.7141999999999
$ws.Range("N85").Value = -3936
$ws.Range("H126").Value = 4120
$ws.Range("I126").Value = 3722.2222
$ws.Range("K126").Value = 11166.6666
$ws.Range("M126").Value = -8696.6666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4833.3335
$ws.Range("H65").Value = 4833.3335
$ws.Range("H81").Value = 2057.7
$ws.Range("I81").Value = 503.66666
$ws.Range("K81").Value = 1007.33332
$ws.Range("M81").Value = 53.66668000000004
$ws.Range("H84").Value = 2057.7
$ws.Range("I84").Value = 503.66666
$ws.Range("K84").Value = 5036.6666
$ws.Range("M84").Value = 267.3334000000004
$ws.Range("H107").Value = 49783856
$ws.Range("I107").Value = 71428800
$ws.Range("J107").Value = 6493966.5
$ws.Range("K107").Value = 214286400
$ws.Range("L107").Value = 19481899.5
$ws.Range("M107").Value = -214284480
$ws.Range("N107").Value = -19485739.5
$ws.Range("H136").Value = 33301984
$ws.Range("I136").Value = 39703670
$ws.Range("J136").Value = 13200
$ws.Range("K136").Value = 119111010
$ws.Range("L136").Value = 39600
$ws.Range("M136").Value = -119108460
$ws.Range("N136").Value = -44700

